# Update driverAthletics Excel file for 2021-2025 with new data:
# - Add two new columns: "Wins" (E) and "Race Starts" (F)
# - Unify formatting across the whole table (Arial, not bold, centered, wrapped)
# - Adjust row 1 height and selection to reflect the new layout

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells -------------------------------------------------
$ws.Range("E1").Value = "Wins"
$ws.Range("F1").Value = "Race Starts"

# --- New data: Wins (E) and Race Starts (F) per driver ----------------
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 163

$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 85

$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 119

$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 212

$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 141

$ws.Range("E7").Value = 61
$ws.Range("F7").Value = 205

$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 45

$ws.Range("E9").Value = 7
$ws.Range("F9").Value = 129

$ws.Range("E10").Value = 8
$ws.Range("F10").Value = 259

$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 22

$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 194

$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 129

$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 67

$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 184

$ws.Range("E16").Value = 104
$ws.Range("F16").Value = 365

$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 6

$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 13

$ws.Range("E19").Value = 10
$ws.Range("F19").Value = 245

$ws.Range("E20").Value = 6
$ws.Range("F20").Value = 274

$ws.Range("E21").Value = 32
$ws.Range("F21").Value = 406

$ws.Range("E22").Value = 1
$ws.Range("F22").Value = 119

$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 88

# --- Unify formatting: the whole table now shares one plain Arial,
#     centered, wrapped style (the previous bold header font is gone). ---
$allData = $ws.Range("A1:F23")
$allData.Font.Name = "Arial"
$allData.Font.Bold = $false
$allData.Font.Size = 11
$allData.Font.Color = 1907739
$allData.HorizontalAlignment = -4108
$allData.VerticalAlignment = -4108
$allData.WrapText = $true
$allData.ReadingOrder = 1

# --- Row 1 grows to match the rest of the header band ------------------
$ws.Rows.Item(1).RowHeight = 28

# --- Selection moves to A2 ----------------------------------------------
$ws.Range("A2").Select()
